$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F:V) between rows 88 and 89 ---
# (columns A:E - index/pais/torneio/temporada/data_partida - stay unchanged
# since both rows share the same match date/competition values)
for ($col = 6; $col -le 22; $col++) {
    $v88 = $ws.Cells.Item(88, $col).Value2
    $v89 = $ws.Cells.Item(89, $col).Value2
    $ws.Cells.Item(88, $col).Value = $v89
    $ws.Cells.Item(89, $col).Value = $v88
}

# --- Append new row 141 (Las Palmas vs Getafe) ---
# Copy cell formatting from row 140 first so the new row matches the
# existing style indices (bold/bordered index column, date-formatted
# data_partida column) instead of minting new styles.
$ws.Cells.Item(140, 1).Copy($ws.Cells.Item(141, 1))
$ws.Cells.Item(140, 5).Copy($ws.Cells.Item(141, 5))

$ws.Cells.Item(141, 1).Value = 140
$ws.Cells.Item(141, 2).Value = "spain"
$ws.Cells.Item(141, 3).Value = "laliga"
$ws.Cells.Item(141, 4).Value = "2023-2024"
$ws.Cells.Item(141, 5).Value = 45261.875
$ws.Cells.Item(141, 6).Value = "Las Palmas"
$ws.Cells.Item(141, 7).Value = 2
$ws.Cells.Item(141, 8).Value = "Getafe"
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 2.35
$ws.Cells.Item(141, 11).Value = "15/11/2023 16:02"
$ws.Cells.Item(141, 12).Value = 2.41
$ws.Cells.Item(141, 13).Value = "01/12/2023 20:53"
$ws.Cells.Item(141, 14).Value = 2.97
$ws.Cells.Item(141, 15).Value = "15/11/2023 16:02"
$ws.Cells.Item(141, 16).Value = 2.98
$ws.Cells.Item(141, 17).Value = "01/12/2023 20:51"
$ws.Cells.Item(141, 18).Value = 3.36
$ws.Cells.Item(141, 19).Value = "15/11/2023 16:02"
$ws.Cells.Item(141, 20).Value = 3.62
$ws.Cells.Item(141, 21).Value = "01/12/2023 20:52"
$ws.Cells.Item(141, 22).Value = "https://www.betexplorer.com/football/spain/laliga/las-palmas-getafe/O8jgIi7r/"
